$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description for "Cerebral Organoids day 40" row (B3) with the new, longer text
$ws.Range("B3").Value = "Cerebral Organoids day 40 compared to hPSC and this is a long wall of text. I write a few more words to make it even longer. This shows that we can have a very verbose description of the data set amd it wil stil be displayed correctly in the app. "

# Update the selected cell in the sheet view to C4
$ws.Range("C4").Select()
